# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the data rows (A:R) of the sheet: each
# destination row ends up containing the values that used to live in a
# different (source) row. Row 25 is untouched. We therefore:
#   1. Snapshot every affected row's full A:R contents first.
#   2. Write the snapshot back out according to the destination<-source
#      mapping, so that rows are effectively permuted without relying on
#      write order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 18   # column R
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,26)

# destination row -> source row (values currently sitting in source row
# become the new values of destination row)
$mapping = @{
    2  = 23
    3  = 18
    4  = 7
    5  = 20
    6  = 4
    7  = 5
    8  = 12
    9  = 2
    10 = 24
    11 = 21
    12 = 19
    13 = 22
    14 = 6
    15 = 14
    16 = 17
    17 = 26
    18 = 11
    19 = 13
    20 = 15
    21 = 10
    22 = 16
    23 = 8
    24 = 9
    26 = 3
}

# 1) Snapshot current values of every row we might read from.
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += ,($ws.Cells.Item($r, $c).Value())
    }
    $snapshot[$r] = $rowVals
}

# 2) Write values back according to the mapping, using the snapshot so
#    that earlier writes never affect later reads.
foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c - 1]
    }
}
